$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{
    2 = @(1, 0.3333333333333333, 0.07899133333333333, 0.236974, 0.08952820636815619, 0.0998891403954695, 3, 1, 0.6598136666666666, 1.979441, 0.007704735356083927, 0.008484678519943686, 0.05211956128155555, 0.469076051534, 0.0006897911369715111, 0.0008475272438890792)
    3 = @(1, 0.3333333333333333, 0.07899133333333333, 0.236974, 0.08952820636815619, 0.0998891403954695, 3, 1, 61.06015933333333, 183.180478, 0.7130079175842846, 0.7851850431306702, 4.823223399285777, 43.409010593572, 0.06383431998761513, 0.07843145900970229)
    4 = @(1, 0.3333333333333333, 0.07899133333333333, 0.236974, 0.08952820636815619, 0.0998891403954695, 2, 0.6666666666666666, 0.1563486666666667, 0.469046, 0.001825704984300993, 0.002010519394650058, 0.01235018964488889, 0.111151706804, 0.0001634520926018706, 0.000200829054080014)
    5 = @(1, 0.3333333333333333, 0.07899133333333333, 0.236974, 0.08952820636815619, 0.0998891403954695, 2, 1, 23.6163295, 47.232659, 0.2757711427815902, 0.2024581319964196, 1.865485355644333, 11.192912133866, 0.02468929578133247, 0.02022336877119485)
    6 = @(1, 0.3333333333333333, 0.07899133333333333, 0.236974, 0.08952820636815619, 0.0998891403954695, 2, 0.6666666666666666, 0.14477, 0.43431, 0.00169049929374041, 0.001861626958316384, 0.01143557532666667, 0.10292017794, 0.0001513473696352137, 0.0001859563166032562)
    7 = @(2, 0.6666666666666666, 0.296947, 0.890841, 0.336557583908845, 0.3755067717093034, 3, 1, 0.6598136666666666, 1.979441, 0.007704735356083927, 0.008484678519943686, 0.1959296888756667, 1.763367199881, 0.002593087116100661, 0.003186054240015323)
    8 = @(2, 0.6666666666666666, 0.296947, 0.890841, 0.336557583908845, 0.3755067717093034, 3, 1, 61.06015933333333, 183.180478, 0.7130079175842846, 0.7851850431306702, 18.13163113355533, 163.184680201998, 0.2399682220500437, 0.2948423007404281)
    9 = @(2, 0.6666666666666666, 0.296947, 0.890841, 0.336557583908845, 0.3755067717093034, 2, 0.6666666666666666, 0.1563486666666667, 0.469046, 0.001825704984300993, 0.002010519394650058, 0.04642726752066667, 0.417845407686, 0.000614454858446678, 0.0007549636473439862)
    10 = @(2, 0.6666666666666666, 0.296947, 0.890841, 0.336557583908845, 0.3755067717093034, 2, 1, 23.6163295, 47.232659, 0.2757711427815902, 0.2024581319964196, 7.0127981960365, 42.07678917621899, 0.0928128695263531, 0.07602439955227154)
    11 = @(2, 0.6666666666666666, 0.296947, 0.890841, 0.336557583908845, 0.3755067717093034, 2, 0.6666666666666666, 0.14477, 0.43431, 0.00169049929374041, 0.001861626958316384, 0.04298901719000001, 0.38690115471, 0.0005689503579008812, 0.0006990535292443954)
    12 = @(1, 0.5, 0.27455, 0.5491, 0.3111729859610415, 0.2314563074056745, 3, 1, 0.6598136666666666, 1.979441, 0.007704735356083927, 0.008484678519943686, 0.1811518421833333, 1.0869110531, 0.002397505506792244, 0.001963832359750409)
    13 = @(1, 0.5, 0.27455, 0.5491, 0.3111729859610415, 0.2314563074056745, 3, 1, 61.06015933333333, 183.180478, 0.7130079175842846, 0.7851850431306702, 16.76406674496667, 100.5844004698, 0.221868802728566, 0.1817360307131902)
    14 = @(1, 0.5, 0.27455, 0.5491, 0.3111729859610415, 0.2314563074056745, 2, 0.6666666666666666, 0.1563486666666667, 0.469046, 0.001825704984300993, 0.002010519394650058, 0.04292552643333333, 0.2575531586, 0.0005681100714488964, 0.0004653473950531945)
    15 = @(1, 0.5, 0.27455, 0.5491, 0.3111729859610415, 0.2314563074056745, 2, 1, 23.6163295, 47.232659, 0.2757711427815902, 0.2024581319964196, 6.483863264225, 25.9354530569, 0.08581252994123613, 0.04686021163614192)
    16 = @(1, 0.5, 0.27455, 0.5491, 0.3111729859610415, 0.2314563074056745, 2, 0.6666666666666666, 0.14477, 0.43431, 0.00169049929374041, 0.001861626958316384, 0.03974660350000001, 0.238479621, 0.0005260377129982351, 0.0004308853015387679)
    17 = @(1, 0.3333333333333333, 0.2318183333333333, 0.695455, 0.2627412237619573, 0.2931477804895526, 3, 1, 0.6598136666666666, 1.979441, 0.007704735356083927, 0.008484678519943686, 0.1529569045172222, 1.376612140655, 0.002024351596219511, 0.002487264676288874)
    18 = @(1, 0.3333333333333333, 0.2318183333333333, 0.695455, 0.2627412237619573, 0.2931477804895526, 3, 1, 61.06015933333333, 183.180478, 0.7130079175842846, 0.7851850431306702, 14.15486436972111, 127.39377932749, 0.1873365728180597, 0.2301752526673496)
    19 = @(1, 0.3333333333333333, 0.2318183333333333, 0.695455, 0.2627412237619573, 0.2931477804895526, 2, 0.6666666666666666, 0.1563486666666667, 0.469046, 0.001825704984300993, 0.002010519394650058, 0.03624448732555556, 0.32620038593, 0.000479687961803548, 0.0005893792981728636)
    20 = @(1, 0.3333333333333333, 0.2318183333333333, 0.695455, 0.2627412237619573, 0.2931477804895526, 2, 1, 23.6163295, 47.232659, 0.2757711427815902, 0.2024581319964196, 5.474698144140834, 32.848188864845, 0.07245644753266846, 0.05935015203681129)
    21 = @(1, 0.3333333333333333, 0.2318183333333333, 0.695455, 0.2627412237619573, 0.2931477804895526, 2, 0.6666666666666666, 0.14477, 0.43431, 0.00169049929374041, 0.001861626958316384, 0.03356034011666667, 0.3020430610500001, 0.0004441638532060798, 0.000545731810929965)
}

foreach ($r in $data.Keys) {
    $rowvals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).Value = $rowvals[$i]
    }
}
